$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (applied as exact text, matching the
# source diff, regardless of whether the text happens to look numeric).
$changes = @{
    'D2' = '30.125.11'
    'E2' = '  -0.54%  '
    'D3' = '1.857.66'
    'E3' = '  -0.63%  '
    'E4' = '  +0.07%  '
    'D5' = '233.70'
    'E5' = '  -0.68%  '
    'E6' = '  +0.02%  '
    'D7' = '0.4676'
    'E7' = '  -0.43%  '
    'B8' = 'Cardano'
    'C8' = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
    'D8' = '0.2830'
    'E8' = '  -1.59%  '
    'B9' = 'Dogecoin'
    'C9' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'D9' = '0.06444'
    'E9' = '  -2.05%  '
    'B10' = 'Solana'
    'C10' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D10' = '20.96'
    'E10' = '  -3.85%  '
    'B11' = 'TRON'
    'C11' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D11' = '0.07754'
    'E11' = '  -3.11%  '
    'B12' = 'WrappedEther'
    'C12' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D12' = '1.862.86'
    'E12' = '  -0.36%  '
    'B13' = 'Litecoin'
    'C13' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D13' = '93.36'
    'E13' = '  -4.03%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D14' = '0.6774'
    'E14' = '  -1.25%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '5.032'
    'E15' = '  -1.97%  '
    'B16' = 'BitcoinCash'
    'C16' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D16' = '266.00'
    'E16' = '  -1.77%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '30.119.83'
    'E17' = '  -0.54%  '
    'B18' = 'Avalanche'
    'C18' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D18' = '13.30'
    'E18' = '  -5.58%  '
    'B19' = 'ShibaInu'
    'C19' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D19' = '0.000007558'
    'E19' = '  -1.73%  '
    'B20' = 'Dai'
    'C20' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D20' = '1.001'
    'E20' = '  +0.04%  '
    'B21' = 'WrappedliquidstakedEther2.0'
    'C21' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D21' = '2.121.03'
    'E21' = '  +0.34%  '
    'B22' = 'BinanceUSD'
    'C22' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D22' = '1.001'
    'E22' = '  -0.02%  '
    'B23' = 'Uniswap'
    'C23' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D23' = '5.121'
    'E23' = '  -3.08%  '
    'B24' = 'Chainlink'
    'C24' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D24' = '6.085'
    'E24' = '  -2.00%  '
    'B25' = 'Cosmos'
    'C25' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D25' = '9.294'
    'E25' = '  -1.96%  '
    'B26' = 'Monero'
    'C26' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D26' = '164.88'
    'E26' = '  -2.03%  '
    'B27' = 'EthereumClassic'
    'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D27' = '18.46'
    'E27' = '  -2.42%  '
    'B28' = 'LidoDAOToken'
    'C28' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D28' = '1.878'
    'E28' = '  -3.67%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D29' = '1.361'
    'E29' = '  -0.79%  '
    'B30' = 'Stellar'
    'C30' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D30' = '0.09887'
    'E30' = '  +0.04%  '
    'B31' = 'PancakeSwap'
    'C31' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D31' = '1.446'
    'E31' = '  -1.28%  '
    'B32' = 'Filecoin'
    'C32' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D32' = '4.193'
    'E32' = '  -3.89%  '
    'B33' = 'InternetComputer(DFINITY)'
    'C33' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D33' = '3.967'
    'E33' = '  -2.49%  '
    'B34' = 'Hedera'
    'C34' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D34' = '0.04640'
    'E34' = '  -1.31%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '1.111'
    'E35' = '  -2.30%  '
    'B36' = 'ImmutableX'
    'C36' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D36' = '0.6856'
    'E36' = '  -2.07%  '
    'B37' = 'HuobiToken'
    'C37' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D37' = '2.714'
    'E37' = '  +0.05%  '
    'B38' = 'VeChain'
    'C38' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D38' = '0.01827'
    'E38' = '  -2.79%  '
    'B39' = 'MXToken'
    'C39' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D39' = '2.730'
    'E39' = '  +3.70%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D40' = '6.251'
    'E40' = '  -0.66%  '
    'B41' = 'Aave'
    'C41' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D41' = '70.67'
    'E41' = '  -2.50%  '
    'B42' = 'PaxDollar'
    'C42' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D42' = '1.000'
    'E42' = '  +0.01%  '
    'B43' = 'TrustWalletToken'
    'C43' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D43' = '0.8327'
    'E43' = '  -1.18%  '
    'B44' = 'RenderToken'
    'C44' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D44' = '1.869'
    'E44' = '  -4.28%  '
    'B45' = 'Quant'
    'C45' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D45' = '101.81'
    'E45' = '  -1.24%  '
    'B46' = 'TheSandbox'
    'C46' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D46' = '0.4036'
    'E46' = '  -3.16%  '
    'B47' = 'EnergySwap'
    'C47' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D47' = '9.076'
    'E47' = '  -1.05%  '
    'B48' = 'Maker'
    'C48' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D48' = '924.07'
    'E48' = '  -0.08%  '
    'B49' = 'Aptos'
    'C49' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D49' = '6.913'
    'E49' = '  -2.03%  '
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D50' = '33.87'
    'E50' = '  -1.63%  '
    'B51' = 'Cronos'
    'C51' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D51' = '0.05552'
    'E51' = '  -2.24%  '
}

foreach ($addr in $changes.Keys) {
    $rng = $ws.Range($addr)
    # Force text storage so Excel does not reinterpret numeric-looking
    # strings (e.g. "233.70", "1.878", "30.125.11") as numbers and strip
    # meaningful trailing zeros / grouping dots.
    $rng.NumberFormat = "@"
    $rng.Value = $changes[$addr]
    # Restore the default (unstyled) cell style so we do not leave behind
    # a lingering custom number-format style that was not in the original.
    $rng.Style = "Normal"
}
